# "Added client report and report page to owner"
#
# The underlying data table on the `datasheet` worksheet only ever had a
# header row (№ / Mount / Profit) plus a single sample row (2018-12 / 700)
# worth of real report data baked in as a template placeholder. This edit
# strips the canned sample rows back out, leaving just the header row, and
# moves the user's selection/active sheet onto the data sheet so the
# workbook opens ready for the next report to be dropped in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three sample data rows (2018-12/700, 2018-11/900, 2018-10/136),
# shifting everything below them up - only the header row (row 3) remains.
$ws.Rows("4:6").Delete() | Out-Null

# Leave the cursor where the next appended report row would go, and make
# sure the data sheet (rather than the chart sheet) is the one in front.
$ws.Range("B6").Select() | Out-Null
$ws.Activate() | Out-Null
